# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.837.21"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.640.77"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "215.96"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.2584"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "0.06439"
$ws.Range("D10").Value = "20.44"
$ws.Range("E10").Value = "  +5.23%  "
$ws.Range("D11").Value = "0.07807"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.647.07"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "1.867.22"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "0.5624"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "63.45"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "25.854.51"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "4.386"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "193.22"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "9.941"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Value = "6.153"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "1.800"
$ws.Range("E25").Value = "  -4.84%  "
$ws.Range("D26").Value = "139.29"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").Value = "0.1233"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "6.835"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").Value = "15.62"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "1.246"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "0.04956"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "3.293"
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("D33").Value = "3.256"
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("D35").Value = "2.388"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").Value = "0.9042"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").Value = "2.571"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "0.5562"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").Value = "1.133.95"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").Value = "0.9973"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8033"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.469"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").Value = "98.94"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "1.777.37"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("D47").Value = "55.63"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").Value = "7.785"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("D50").Value = "0.05036"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -0.03%  "
